# Applies the "add rotate / change folder structure / densenet201" commit to NOTES.xlsx
# (fixed preprocess.py (add rotate) change folder structure change net to densenet201)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel constant values (avoid depending on enum types being registered)
$xlCenter = -4108
$xlLeft   = -4131

# ---------------------------------------------------------------------------
# 1. Fill in the two new experiment rows (6-9) with their submission name,
#    model name, date, score and note.  (rows 6-9 already had a date in C
#    and row 6 already had a note; everything else is new.)
# ---------------------------------------------------------------------------

$ws.Range("A6").Value2 = "novoting_full_densenet121_AutoWtdCE_2020-12-03_1-12_epoch49"
$ws.Range("B6").Value2 = "full_densenet121_AutoWtdCE_2020-12-03_1-12_epoch49"
$ws.Range("C6").Value2 = 44168
$ws.Range("D6").Value2 = 0.416
$ws.Range("E6").Value2 = "add random increase brightness to preprocess"

$ws.Range("A7").Value2 = "voting_full_densenet121_AutoWtdCE_2020-12-03_1-12_epoch49"
$ws.Range("B7").Value2 = "full_densenet121_AutoWtdCE_2020-12-03_1-12_epoch49"
$ws.Range("C7").Value2 = 44168
$ws.Range("D7").Value2 = 0.16
$ws.Range("E7").Value2 = "add random increase brightness to preprocess and using voting"

$ws.Range("A8").Value2 = "novoting_full_densenet121_AutoWtdCE_2020-12-03_18-35_epoch48"
$ws.Range("B8").Value2 = "full_densenet121_AutoWtdCE_2020-12-03_18-35_epoch49"
$ws.Range("C8").Value2 = 44168
$ws.Range("D8").Value2 = 0.431
$ws.Range("E8").Value2 = "add random increase brightness to preprocess"

$ws.Range("A9").Value2 = "voting_full_densenet121_AutoWtdCE_2020-12-03_18-35_epoch49"
$ws.Range("B9").Value2 = "full_densenet121_AutoWtdCE_2020-12-03_18-35_epoch49"
$ws.Range("C9").Value2 = 44168
$ws.Range("D9").Value2 = 0.165
$ws.Range("E9").Value2 = "add random increase brightness to preprocess and using voting"

# ---------------------------------------------------------------------------
# 2. Re-align the whole data body of the table (rows 2-49) the way the
#    workbook author did: Submissions/Model centred vertically, Date and
#    Score centred both ways, Note centred vertically (keeping its wrap).
#
#    A scratch cell (Z1) is used to build up each combination of alignment
#    settings exactly once and then the resulting format is copied onto the
#    target range with PasteSpecial(xlPasteFormats) - this avoids Excel
#    recording one intermediate style per property assignment.
# ---------------------------------------------------------------------------

$xlPasteFormats = -4122

# Submissions / Model columns -> vertical-center only
$ws.Range("Z1").VerticalAlignment = $xlCenter
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("A2:B49").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# Note column -> vertical-center, keep the existing word-wrap
$ws.Range("Z1").VerticalAlignment = $xlCenter
$ws.Range("Z1").WrapText = $true
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("E2:E49").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# Score column -> horizontal+vertical center
$ws.Range("Z1").HorizontalAlignment = $xlCenter
$ws.Range("Z1").VerticalAlignment = $xlCenter
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("D2:D9").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# Date column -> keep its existing date number format, add horizontal+vertical center
$ws.Range("C2").HorizontalAlignment = $xlCenter
$ws.Range("C2").VerticalAlignment = $xlCenter
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C2:C49").PasteSpecial($xlPasteFormats) | Out-Null

# ---------------------------------------------------------------------------
# 3. Column widths - Model column got narrower once the new (shorter)
#    model-name column was reformatted by the author.
# ---------------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth = 53.5703125

# ---------------------------------------------------------------------------
# 4. Selection left by the author at the end of editing.
# ---------------------------------------------------------------------------

[void]$ws.Range("B14").Select()
